# Adds two new columns, I (I0) and J (IF), to the sheet, mirroring the
# existing header style used by B1:H1, and fills in data rows 2-72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns I1 = "I0", J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the bold/centered/bordered style used by the other header cells
# (B1:H1) by copying their formatting onto the new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows 2-72: column I ("I0") values ---
$iValues = @(8,9,9,5,5,9,8,4,13,9,6,7,5,6,7,6,6,6,8,6,8,6,8,8,7,7,7,9,8,9,7,8,7,7,9,9,9,9,9,9,8,8,9,9,9,9,8,9,10,8,8,9,8,7,8,8,8,8,7,7,8,8,7,8,7,6,7,8,4,5,3)

# --- Data rows 2-72: column J ("IF") values ---
$jValues = @(8,9,9,5,5,9,8,4,13,9,6,8,5,6,7,6,6,6,8,7,8,6,8,9,7,7,7,9,8,9,7,8,7,7,9,9,9,9,9,9,8,9,9,9,9,9,8,9,10,8,8,9,8,8,8,8,8,8,7,8,8,8,7,9,7,6,7,8,4,5,3)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
